$wb = $excel.ActiveWorkbook

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4600
$ws.Range("I62").Value = 4600
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4600
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3976
$ws.Range("N62").ClearContents()

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 4600
$ws.Range("I65").Value = 4600
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 23000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -19880
$ws.Range("N65").ClearContents()

# ALC row 124
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H124").Value = 44900
$ws.Range("J124").Value = 44900
$ws.Range("L124").Value = 44900
$ws.Range("N124").Value = -54720

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3263.717
$ws.Range("I132").Value = 3706.0264
$ws.Range("J132").Value = 2143.2
$ws.Range("K132").Value = 11118.0792
$ws.Range("L132").Value = 6429.599999999999
$ws.Range("M132").Value = -8588.0792
$ws.Range("N132").Value = -11489.6

# ARM row 7
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 152000
$ws.Range("J7").Value = 152000
$ws.Range("L7").Value = 152000
$ws.Range("N7").Value = -152228

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3206.689
$ws.Range("I61").Value = 2891.838
$ws.Range("K61").Value = 2891.838
$ws.Range("M61").Value = -2679.838

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1460.5312
$ws.Range("I74").Value = 1384.1
$ws.Range("K74").Value = 1384.1
$ws.Range("M74").Value = -510.0999999999999

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1460.5312
$ws.Range("I77").Value = 1384.1
$ws.Range("K77").Value = 6920.5
$ws.Range("M77").Value = -2552.5

# ARM row 128
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H128").Value = 99784.5
$ws.Range("J128").Value = 99784.5
$ws.Range("L128").Value = 99784.5
$ws.Range("N128").Value = -109744.5

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2368.7144
$ws.Range("I132").Value = 1531.4423
$ws.Range("K132").Value = 4594.3269
$ws.Range("M132").Value = -2064.3269

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3206.689
$ws.Range("I136").Value = 2891.838
$ws.Range("K136").Value = 8675.514000000001
$ws.Range("M136").Value = -6125.514000000001

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11567122
$ws.Range("I134").Value = 2464065.5
$ws.Range("J134").Value = 55565230
$ws.Range("K134").Value = 7392196.5
$ws.Range("L134").Value = 166695690
$ws.Range("M134").Value = -7389661.5
$ws.Range("N134").Value = -166700760

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 178.2
$ws.Range("I7").Value = 96.85714
$ws.Range("J7").Value = 368
$ws.Range("K7").Value = 96.85714
$ws.Range("L7").Value = 368
$ws.Range("M7").Value = 16.14286
$ws.Range("N7").Value = -594

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3816.5862
$ws.Range("I31").Value = 2210.4
$ws.Range("K31").Value = 2210.4
$ws.Range("M31").Value = -1915.4

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3816.5862
$ws.Range("I34").Value = 2210.4
$ws.Range("K34").Value = 2210.4
$ws.Range("M34").Value = -2008.4

# CRP row 75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 90420.14
$ws.Range("J75").Value = 102788.6
$ws.Range("L75").Value = 102788.6
$ws.Range("N75").Value = -104784.6

# CRP row 78
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H78").Value = 90420.14
$ws.Range("J78").Value = 102788.6
$ws.Range("L78").Value = 308365.8
$ws.Range("N78").Value = -318349.8

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3490.8333
$ws.Range("I86").Value = 3490.8333
$ws.Range("K86").Value = 3490.8333
$ws.Range("M86").Value = -2367.8333

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3490.8333
$ws.Range("I89").Value = 3490.8333
$ws.Range("K89").Value = 17454.1665
$ws.Range("M89").Value = -11838.1665

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1676
$ws.Range("I132").Value = 1569
$ws.Range("K132").Value = 4707
$ws.Range("M132").Value = -2177

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1418.4445
$ws.Range("I17").Value = 3394.6667
$ws.Range("J17").Value = 430.33334
$ws.Range("K17").Value = 10184.0001
$ws.Range("L17").Value = 1291.00002
$ws.Range("M17").Value = -10015.0001
$ws.Range("N17").Value = -1629.00002

# CUL row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 2415.625
$ws.Range("I18").Value = 65.2
$ws.Range("K18").Value = 195.6
$ws.Range("M18").Value = -26.60000000000002

# CUL row 20
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 2929.6667
$ws.Range("J20").Value = 3990
$ws.Range("L20").Value = 11970
$ws.Range("N20").Value = -12424

# CUL row 21
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()

# CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 128
$ws.Range("I26").Value = 128
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 384
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -96
$ws.Range("N26").ClearContents()

# CUL row 32
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 555
$ws.Range("I34").Value = 140.9
$ws.Range("J34").Value = 1935.3334
$ws.Range("K34").Value = 422.7
$ws.Range("L34").Value = 5806.0002
$ws.Range("M34").Value = -338.7
$ws.Range("N34").Value = -5974.0002

# CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3872.5833
$ws.Range("J39").Value = 4441.222
$ws.Range("L39").Value = 13323.666
$ws.Range("N39").Value = -13911.666

# CUL row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 3125.7144
$ws.Range("J40").Value = 5412.5
$ws.Range("L40").Value = 21650
$ws.Range("N40").Value = -21788

# CUL row 51
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 3874
$ws.Range("I51").Value = 444
$ws.Range("J51").Value = 4560
$ws.Range("K51").Value = 1332
$ws.Range("L51").Value = 13680
$ws.Range("M51").Value = -872
$ws.Range("N51").Value = -14600

# CUL row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1260.25
$ws.Range("I55").Value = 1084.5714
$ws.Range("J55").Value = 2490
$ws.Range("K55").Value = 3253.7142
$ws.Range("L55").Value = 7470
$ws.Range("M55").Value = -3076.7142
$ws.Range("N55").Value = -7824

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2747.5789
$ws.Range("I102").Value = 2646.4614
$ws.Range("K102").Value = 2646.4614
$ws.Range("M102").Value = -1024.4614

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1431.7727
$ws.Range("I132").Value = 1219
$ws.Range("J132").Value = 2070.0908
$ws.Range("K132").Value = 3657
$ws.Range("L132").Value = 6210.2724
$ws.Range("M132").Value = -1127
$ws.Range("N132").Value = -11270.2724

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2828
$ws.Range("J7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("N7").Value = -3224

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4790.56
$ws.Range("I40").Value = 4925.263
$ws.Range("J40").Value = 4364
$ws.Range("K40").Value = 4925.263
$ws.Range("L40").Value = 4364
$ws.Range("M40").Value = -4789.263
$ws.Range("N40").Value = -4636

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3282.2
$ws.Range("I82").Value = 2662.3333
$ws.Range("K82").Value = 2662.3333
$ws.Range("M82").Value = -2301.3333

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3282.2
$ws.Range("I85").Value = 2662.3333
$ws.Range("K85").Value = 2662.3333
$ws.Range("M85").Value = -1414.3333

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4168.3335

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2828
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2215.282
$ws.Range("I132").Value = 1709.3704
$ws.Range("J132").Value = 3353.5833
$ws.Range("K132").Value = 5128.1112
$ws.Range("L132").Value = 10060.7499
$ws.Range("M132").Value = -2598.1112
$ws.Range("N132").Value = -15120.7499

# WVR row 130
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 122971.5
$ws.Range("J130").Value = 122971.5
$ws.Range("L130").Value = 122971.5
$ws.Range("N130").Value = -133011.5
